$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1: i / Xi source data (column D, rows 4-12) ---
$ws.Range("D4").Value = 9
$ws.Range("D5").Value = 11
$ws.Range("D6").Value = 14
$ws.Range("D8").Value = 19
$ws.Range("D9").Value = 22
$ws.Range("D10").Value = 27
$ws.Range("D11").Value = 33
$ws.Range("D12").Value = 39

# --- Model of Transition Probabilities (rows 29-35) ---
$ws.Range("E29").Value = 0.08
$ws.Range("F29").Value = 0.11

$ws.Range("E30").Value = 0.4
$ws.Range("F30").Value = 0.3
$ws.Range("G30").Value = 0.9

$ws.Range("G31").Value = 0.8

$ws.Range("E33").Value = 0.3
$ws.Range("E34").Value = 0.9
$ws.Range("E35").Value = 0.5

# --- Korkoran model table (rows 49-55) ---
$ws.Range("D49").Value = 0.066
$ws.Range("E49").Value = 7

$ws.Range("D50").Value = 0.3
$ws.Range("E50").Value = 20

$ws.Range("D51").Value = 0.08
$ws.Range("E51").Value = 8

$ws.Range("D52").Value = 0.3

$ws.Range("D53").Value = 0.2
$ws.Range("E53").Value = 13

$ws.Range("D54").Value = 0.0763

$ws.Range("D55").Value = 0.2
$ws.Range("E55").Value = 6

# --- Highlight J49 / J50 with a yellow fill (new cells picking up new styles) ---
$ws.Range("J50").Interior.Color = 65535

$ws.Range("F3").Copy()
$ws.Range("J49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("J49").Interior.Color = 65535

# --- View state: zoom in, scroll down, move the active selection ---
[void]$excel.Goto($ws.Range("A25"), $true)
$ws.Range("E36").Select() | Out-Null
$excel.ActiveWindow.Zoom = 115

Write-Host "Edit applied"
